$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original text formatting
# (Excel would otherwise auto-convert numeric-looking strings to numbers,
# dropping significant trailing/leading zeros and thousands-dot formatting).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.100.12'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.891.82'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.50'
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5185'
$ws.Range("E7").Value = '  +2.91%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3755'
$ws.Range("E8").Value = '  +3.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07217'
$ws.Range("E9").Value = '  +0.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.11'
$ws.Range("E10").Value = '  +2.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9014'
$ws.Range("E11").Value = '  +1.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07658'
$ws.Range("E12").Value = '  +2.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.915.10'
$ws.Range("E13").Value = '  +2.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.35'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.237'
$ws.Range("E15").Value = '  +0.29%  '

$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008513'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.43'
$ws.Range("E18").Value = '  +1.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.144.76'
$ws.Range("E20").Value = '  +0.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.060'
$ws.Range("E21").Value = '  +0.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.128.83'
$ws.Range("E22").Value = '  +1.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.59'
$ws.Range("E23").Value = '  +2.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.382'
$ws.Range("E24").Value = '  -0.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.300'
$ws.Range("E25").Value = '  +11.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '145.36'
$ws.Range("E26").Value = '  -1.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.06'
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.731'
$ws.Range("E28").Value = '  -2.89%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.36'
$ws.Range("E29").Value = '  +1.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.978'
$ws.Range("E30").Value = '  +6.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.793'
$ws.Range("E31").Value = '  +2.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09210'
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05054'
$ws.Range("E33").Value = '  -1.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.240'
$ws.Range("E34").Value = '  +7.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7761'
$ws.Range("E35").Value = '  +4.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.978'
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.288'
$ws.Range("E37").Value = '  +3.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.595'
$ws.Range("E38").Value = '  +1.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5649'
$ws.Range("E39").Value = '  +1.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01988'
$ws.Range("E40").Value = '  -0.39%  '

$ws.Range("E41").Value = '  +0.36%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.037'
$ws.Range("E42").Value = '  +6.01%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.63'
$ws.Range("E43").Value = '  +3.39%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.642'
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("E45").Value = '  +2.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4845'
$ws.Range("E46").Value = '  +3.48%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.20'
$ws.Range("E47").Value = '  +1.87%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9996'
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.595'
$ws.Range("E49").Value = '  +2.58%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.68'
$ws.Range("E50").Value = '  +2.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.95'
